$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new user row (email + password) below the existing data.
$ws.Range("A7").Value = "aba123@gmail.com"
$ws.Range("B7").Value = "322aas"

# Make the new email a real mailto hyperlink, matching the existing rows.
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:aba123@gmail.com") | Out-Null

$ws.Range("A14").Select()
